$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values
$ws.Range("C6").Value = 21083
$ws.Range("C12").Value = 6793

# Update active cell selection (cosmetic change seen in diff)
$ws.Range("D22").Select()
